$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 21 (entry 8) ----
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = 45302
$ws.Cells.Item(21, 3).Value = "IPS/1803"
$ws.Cells.Item(21, 4).Value = "Vissu Virgincar & Sons"
$ws.Cells.Item(21, 5).Value = 40500
$ws.Cells.Item(21, 6).Formula = "=E21"

# ---- Row 23 (entry 9) ----
$ws.Cells.Item(23, 1).Value = 9
$ws.Cells.Item(23, 2).Value = 45306
$ws.Cells.Item(23, 3).Value = 1436
$ws.Cells.Item(23, 4).Value = "Aquachemitech"
$ws.Cells.Item(23, 5).Value = 26400
$ws.Cells.Item(23, 6).Formula = "=E23"

# ---- Row 25 (entry 10) ----
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = 45308
$ws.Cells.Item(25, 3).Value = "SLH/3633"
$ws.Cells.Item(25, 4).Value = "Shree Laxmi Lighting Hub"
$ws.Cells.Item(25, 5).Value = 2877
$ws.Cells.Item(25, 6).Formula = "=E25"

# Apply the same look-and-feel used by the other numbered entries (row 6
# carries the same A:F style pattern as the new rows) by copying formats only.
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A21:F21").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:F23").PasteSpecial(-4122) | Out-Null
$ws.Range("A25:F25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(21).RowHeight = 14.4
$ws.Rows.Item(23).RowHeight = 14.4
$ws.Rows.Item(25).RowHeight = 14.4

$ws.Range("G19").Select()
